$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E price/volume cells to keep their literal text
# representation (these columns store pre-formatted display strings, e.g.
# "596.20" or "2.981.17", which must not be re-interpreted as numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '63.044.04'
$ws.Range("E2").Value = '  +3.37%  '

$ws.Range("D3").Value = '2.984.17'
$ws.Range("E3").Value = '  +2.38%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '596.20'
$ws.Range("E5").Value = '  +1.05%  '

$ws.Range("D6").Value = '146.43'
$ws.Range("E6").Value = '  +0.97%  '

$ws.Range("D8").Value = '2.983.04'
$ws.Range("E8").Value = '  +2.48%  '

$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("D10").Value = '7.45'
$ws.Range("E10").Value = '  +7.21%  '

$ws.Range("D11").Value = '0.145'
$ws.Range("E11").Value = '  +2.57%  '

$ws.Range("E12").Value = '  +3.23%  '

$ws.Range("D13").Value = '0.0000237'
$ws.Range("E13").Value = '  +5.59%  '

$ws.Range("D14").Value = '33.73'
$ws.Range("E14").Value = '  +0.80%  '

$ws.Range("E15").Value = '  +0.34%  '

$ws.Range("D16").Value = '3.475.31'
$ws.Range("E16").Value = '  +2.38%  '

$ws.Range("D17").Value = '62.898.86'
$ws.Range("E17").Value = '  +3.40%  '

$ws.Range("D18").Value = '6.78'
$ws.Range("E18").Value = '  +1.51%  '

$ws.Range("D19").Value = '2.982.93'
$ws.Range("E19").Value = '  +2.36%  '

$ws.Range("D20").Value = '444.96'
$ws.Range("E20").Value = '  +2.39%  '

$ws.Range("D21").Value = '13.62'
$ws.Range("E21").Value = '  +1.85%  '

$ws.Range("D22").Value = '0.679'
$ws.Range("E22").Value = '  +0.58%  '

$ws.Range("D23").Value = '7.18'
$ws.Range("E23").Value = '  +1.12%  '

$ws.Range("D24").Value = '82.38'
$ws.Range("E24").Value = '  +1.08%  '

$ws.Range("E25").Value = '  +1.70%  '

$ws.Range("D26").Value = '12.15'
$ws.Range("E26").Value = '  +3.16%  '

$ws.Range("E27").Value = '  -1.20%  '

$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("E29").Value = '  +1.66%  '

$ws.Range("D30").Value = '7.18'
$ws.Range("E30").Value = '  +2.82%  '

$ws.Range("D31").Value = '2.15'
$ws.Range("E31").Value = '  -5.91%  '

$ws.Range("D32").Value = '26.70'
$ws.Range("E32").Value = '  +0.78%  '

$ws.Range("E33").Value = '  -0.15%  '

$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").Value = '0.0₃0890'
$ws.Range("E35").Value = '  +2.43%  '

$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -1.40%  '

$ws.Range("D37").Value = '5.70'
$ws.Range("E37").Value = '  +1.51%  '

$ws.Range("E38").Value = '  +3.95%  '

$ws.Range("E39").Value = '  +0.41%  '

$ws.Range("D40").Value = '2.97'
$ws.Range("E40").Value = '  -1.00%  '

$ws.Range("D41").Value = '8.66'
$ws.Range("E41").Value = '  +0.98%  '

$ws.Range("E42").Value = '  -2.56%  '

$ws.Range("E43").Value = '  -0.29%  '

$ws.Range("D44").Value = '39.29'
$ws.Range("E44").Value = '  -6.43%  '

$ws.Range("D45").Value = '374.66'
$ws.Range("E45").Value = '  -0.62%  '

$ws.Range("D46").Value = '2.716.21'
$ws.Range("E46").Value = '  +1.07%  '

$ws.Range("E47").Value = '  -0.20%  '

$ws.Range("D48").Value = '135.10'
$ws.Range("E48").Value = '  +1.70%  '

$ws.Range("D50").Value = '23.50'
$ws.Range("E50").Value = '  -1.42%  '

$ws.Range("D51").Value = '0.106'
$ws.Range("E51").Value = '  +0.07%  '

